# Applies the cryptos.xlsx price/volume/coin-listing refresh described in the commit
# "Updated cryptos list on Wed Jun  5 23:48:13 UTC 2024 with GitHub Actions".
#
# The sheet is a straight data dump: column A is a 0-based rank index (untouched),
# B/C/D/E are Coin / Link / Price / Volume(1h). Every row from 2-51 gets a refreshed
# Price and/or Volume, rows 29-51 also shift down by one slot to make room for the new
# "WrappedeETH" listing (each row now shows the coin that used to be one row below it).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '71.081.48'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.75%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.866.52'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.41%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '698.14'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +2.02%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '173.55'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +1.23%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.866.29'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +1.42%  '

$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -0.05%  '

$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +1.19%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.18'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -1.50%  '

$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.08%  '

$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +5.19%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '36.49'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +1.23%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.517.34'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +1.47%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.862.38'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +1.46%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '71.133.09'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.75%  '

$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.18%  '

$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +0.96%  '

$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.07%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.17'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -1.87%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '495.85'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +4.22%  '

$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +1.47%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '85.23'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +2.00%  '

$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +1.71%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '10.75'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +3.98%  '

$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +0.29%  '

$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +1.69%  '

$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.006.21'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +1.08%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.19'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +8.19%  '

$ws.Range('B31').Value = 'Dai'
$ws.Range('C31').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +0.02%  '

$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.67'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +3.53%  '

$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.29'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.69%  '

$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '29.83'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +0.37%  '

$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.178'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -0.75%  '

$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '9.33'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +1.96%  '

$ws.Range('B37').Value = 'RenzoRestakedETH'
$ws.Range('C37').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.817.52'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +1.48%  '

$ws.Range('B38').Value = 'Binance-PegBSC-USD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.05%  '

$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.105'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +2.21%  '

$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.41'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +12.20%  '

$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.08'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +2.20%  '

$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.40'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +0.10%  '

$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.04'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +7.48%  '

$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -0.04%  '

$ws.Range('B45').Value = 'USDe'
$ws.Range('C45').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +0.04%  '

$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '162.85'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +1.86%  '

$ws.Range('B47').Value = 'FLOKI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.000309'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +4.81%  '

$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '48.62'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.75%  '

$ws.Range('B49').Value = 'Arweave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '44.35'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -4.35%  '

$ws.Range('B50').Value = 'TheGraph'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.304'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +1.60%  '

$ws.Range('B51').Value = 'Bittensor'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '418.77'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +5.13%  '
